$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J3 and K3 values
$ws.Range("J3").Value = 1983.102294921875
$ws.Range("K3").Value = 278.021787060975

# Remove J14/K14 and J16/K16 (entire cell contents cleared)
$ws.Range("J14").ClearContents()
$ws.Range("K14").ClearContents()
$ws.Range("J16").ClearContents()
$ws.Range("K16").ClearContents()

# Adjust column widths to reflect new best-fit sizes (closest value achievable
# through the ColumnWidth COM setter, which quantizes to a pixel grid)
$ws.Columns.Item(1).ColumnWidth = 21.8333333333333
$ws.Columns.Item(2).ColumnWidth = 8.8333333333333
$ws.Columns.Item(3).ColumnWidth = 8.8333333333333
$ws.Columns.Item(4).ColumnWidth = 6.5
$ws.Columns.Item(5).ColumnWidth = 6.5
$ws.Columns.Item(6).ColumnWidth = 6.5
$ws.Columns.Item(7).ColumnWidth = 7.6666666666667
$ws.Columns.Item(8).ColumnWidth = 8.8333333333333
$ws.Columns.Item(10).ColumnWidth = 16.6666666666667
